$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(90, 183957, 24648, 51600),
    @(91, 187327, 25085, 54543),
    @(92, 189973, 25549, 57576),
    @(93, 192994, 25969, 60498),
    @(94, 195351, 26384, 63120),
    @(95, 197675, 26644, 64928),
    @(96, 199414, 26977, 66624)
)

$startRow = 92
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
